# Ichthyosis.xlsx edit: refresh "data" sheet export timestamps and add a new
# "metadata" tab (panel/provenance info) right after the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Refresh the per-gene "time_taken" timestamps on the "data" sheet -------
$dataSheet.Range("F2").Value = "2021-10-05 14:34:07.890267"
$dataSheet.Range("F3").Value = "2021-10-05 14:34:07.890276"
$dataSheet.Range("F4").Value = "2021-10-05 14:34:07.890279"
$dataSheet.Range("F5").Value = "2021-10-05 14:34:07.890282"
$dataSheet.Range("F6").Value = "2021-10-05 14:34:07.890285"
$dataSheet.Range("F7").Value = "2021-10-05 14:34:07.890288"
$dataSheet.Range("F8").Value = "2021-10-05 14:34:07.890291"
$dataSheet.Range("F9").Value = "2021-10-05 14:34:07.890293"
$dataSheet.Range("F10").Value = "2021-10-05 14:34:07.890296"
$dataSheet.Range("F11").Value = "2021-10-05 14:34:07.890299"
$dataSheet.Range("F12").Value = "2021-10-05 14:34:07.890302"
$dataSheet.Range("F13").Value = "2021-10-05 14:34:07.890305"
$dataSheet.Range("F14").Value = "2021-10-05 14:34:07.890307"
$dataSheet.Range("F15").Value = "2021-10-05 14:34:07.890311"
$dataSheet.Range("F16").Value = "2021-10-05 14:34:07.890314"
$dataSheet.Range("F17").Value = "2021-10-05 14:34:07.890317"
$dataSheet.Range("F18").Value = "2021-10-05 14:34:07.890321"
$dataSheet.Range("F19").Value = "2021-10-05 14:34:07.890325"
$dataSheet.Range("F20").Value = "2021-10-05 14:34:07.890330"
$dataSheet.Range("F21").Value = "2021-10-05 14:34:07.890334"
$dataSheet.Range("F22").Value = "2021-10-05 14:34:07.890338"
$dataSheet.Range("F23").Value = "2021-10-05 14:34:07.890342"
$dataSheet.Range("F24").Value = "2021-10-05 14:34:07.890347"
$dataSheet.Range("F25").Value = "2021-10-05 14:34:07.890351"
$dataSheet.Range("F26").Value = "2021-10-05 14:34:07.890356"
$dataSheet.Range("F27").Value = "2021-10-05 14:34:07.890361"
$dataSheet.Range("F28").Value = "2021-10-05 14:34:07.890365"
$dataSheet.Range("F29").Value = "2021-10-05 14:34:07.890370"
$dataSheet.Range("F30").Value = "2021-10-05 14:34:07.890374"
$dataSheet.Range("F31").Value = "2021-10-05 14:34:07.890379"
$dataSheet.Range("F32").Value = "2021-10-05 14:34:07.890383"
$dataSheet.Range("F33").Value = "2021-10-05 14:34:07.890388"
$dataSheet.Range("F34").Value = "2021-10-05 14:34:07.890392"
$dataSheet.Range("F35").Value = "2021-10-05 14:34:07.890395"
$dataSheet.Range("F36").Value = "2021-10-05 14:34:07.890398"
$dataSheet.Range("F37").Value = "2021-10-05 14:34:07.890401"
$dataSheet.Range("F38").Value = "2021-10-05 14:34:07.890404"
$dataSheet.Range("F39").Value = "2021-10-05 14:34:07.890406"
$dataSheet.Range("F40").Value = "2021-10-05 14:34:07.890409"
$dataSheet.Range("F41").Value = "2021-10-05 14:34:07.890412"
$dataSheet.Range("F42").Value = "2021-10-05 14:34:07.890415"
$dataSheet.Range("F43").Value = "2021-10-05 14:34:07.890418"
$dataSheet.Range("F44").Value = "2021-10-05 14:34:07.890422"
$dataSheet.Range("F45").Value = "2021-10-05 14:34:07.890424"
$dataSheet.Range("F46").Value = "2021-10-05 14:34:07.890427"
$dataSheet.Range("F47").Value = "2021-10-05 14:34:07.890430"
$dataSheet.Range("F48").Value = "2021-10-05 14:34:07.890433"
$dataSheet.Range("F49").Value = "2021-10-05 14:34:07.890436"
$dataSheet.Range("F50").Value = "2021-10-05 14:34:07.890439"
$dataSheet.Range("F51").Value = "2021-10-05 14:34:07.890442"
$dataSheet.Range("F52").Value = "2021-10-05 14:34:07.890445"
$dataSheet.Range("F53").Value = "2021-10-05 14:34:07.890447"
$dataSheet.Range("F54").Value = "2021-10-05 14:34:07.890450"

# --- Add the new "metadata" sheet, placed right after "data" ---------------
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match the page margins used by the rest of the workbook (0.75/0.75/1/1/.5/.5 in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Copy the header row's style (bold, bordered, centered) from "data" B1:F1
# onto "metadata" B1:F1, then stamp the extra G1 header with the same style.
$dataSheet.Range("B1:F1").Copy($ws.Range("B1:F1"))
$dataSheet.Range("B1").Copy($ws.Range("G1"))
# Copy the "index" column style (data!A2) onto metadata!A2
$dataSheet.Range("A2").Copy($ws.Range("A2"))

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Ichthyosis"
$ws.Range("C2").Value = 124

# data_version "1.1" must stay text, not become the number 1.1
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.1"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "2021-03-30T22:57:37.351967Z"
$ws.Range("F2").Value = "2021-10-05 14:34:07.886749"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/124/?format=json"

# Leave "data" as the selected/active sheet, same as before the edit.
$dataSheet.Select()

Write-Host "metadata sheet added; data sheet timestamps refreshed"
